# fixed crash issue and full model NULL results issue
#
# Rows 66-70 (NTP 2019 28 day study / PFOS / Rat block) had their Sex,
# Model_Type and Serum_Concentration_mg_L values accidentally swapped with
# rows 71-75 - this produced "Male" rows tagged with a bogus "PBPK" model
# type (causing the full PBPK model to crash / return NULL results) while
# the real "Male"+"PBPK" results were sitting on rows 71-75 tagged as
# "Female"+"single-compartment". Swap the F, G and L columns back between
# the two 5-row blocks so each row has the correct Sex / Model_Type /
# Serum_Concentration combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$topRows = 66..70
$bottomRows = 71..75

for ($i = 0; $i -lt 5; $i++) {
    $r1 = $topRows[$i]
    $r2 = $bottomRows[$i]

    $f1 = $ws.Range("F$r1").Value2
    $g1 = $ws.Range("G$r1").Value2
    $l1 = $ws.Range("L$r1").Value2

    $f2 = $ws.Range("F$r2").Value2
    $g2 = $ws.Range("G$r2").Value2
    $l2 = $ws.Range("L$r2").Value2

    $ws.Range("F$r1").Value = $f2
    $ws.Range("G$r1").Value = $g2
    $ws.Range("L$r1").Value = $l2

    $ws.Range("F$r2").Value = $f1
    $ws.Range("G$r2").Value = $g1
    $ws.Range("L$r2").Value = $l1
}

# The table's persisted sort (by PFAS column D) is now stale once the data
# was hand corrected above, so re-apply it against the table's own PFAS
# column (which is already in the correct ascending order) purely to clear
# the old sort-state bookkeeping without actually moving any rows.
$lo = $ws.ListObjects.Item(1)
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("D2:D75"))
$lo.Sort.Apply()

# Flag the now-duplicated PFOS / Rat combination in the corrected block with
# the standard "Duplicate Values" conditional formatting rule.
$rng = $ws.Range("D70:E75")
$fc = $rng.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

# Leave the cursor where the analyst last worked while verifying the fix.
$ws.Range("F64").Select()

$wb.Save()
